$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for column G (K) per row, regenerated from K instead of Strike#
$gValues = @{
    2 = 0
    3 = 0
    4 = 2
    5 = 1
    6 = 1
    7 = 0
    8 = 0
    9 = 0
    10 = 1
    11 = 0
    12 = 3
    13 = 0
    14 = 0
    15 = 1
    17 = 0
    18 = 1
    19 = 2
    20 = 5
    21 = 1
    22 = 1
    23 = 1
    24 = 1
    25 = 1
    26 = 1
    27 = 1
    28 = 1
    29 = 2
    30 = 1
    31 = 4
    32 = 0
    33 = 1
    34 = 2
    35 = 1
    36 = 0
    37 = 2
    38 = 2
    40 = 3
    41 = 3
    42 = 0
    43 = 2
    44 = 2
    45 = 1
    47 = 1
    48 = 0
    49 = 1
    50 = 3
    51 = 1
    52 = 3
    53 = 5
    54 = 3
    55 = 0
    56 = 1
    57 = 0
}

foreach ($row in $gValues.Keys) {
    $ws.Cells.Item([int]$row, 7).Value = $gValues[$row]
}
